# Sara Alert "Invalid Monitorees" template: add three new Race columns
# (Unknown / Other / Refused to Answer) to the header row, right after the
# existing "Sexual Orientation" column, and leave the selection on the
# newly-added range, matching what a user would see after typing the new
# headers in by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing header row runs from A1 through CU1 ("Sexual Orientation" is the
# last column). Append the three new header cells immediately after it.
$ws.Range("CV1").Value = "Race Unknown"
$ws.Range("CW1").Value = "Race Other"
$ws.Range("CX1").Value = "Race Refused to Answer"

# Leave the active selection on the last cell just entered, as the workbook
# was left after the edit.
$ws.Range("CX6").Select()
